$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.177.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.795.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '117.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.12%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '343.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.543'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.93%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0870'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.229.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.800.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.889'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.980.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.99'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0990'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '278.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('E29').Value = '  +3.45%  '
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0825'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.10%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.03'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  +7.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +27.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0376'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +15.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '128.08'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.19%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.26%  '
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.110.91'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.45%  '
$ws.Range('E48').Value = '  +3.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.919'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +21.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.37%  '

Write-Host "Applied cryptos update"
